$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 3).Value = 14
$ws.Cells.Item(5, 3).Value = 24
$ws.Cells.Item(7, 3).Value = 10
$ws.Cells.Item(8, 3).Value = 13
$ws.Cells.Item(11, 3).Value = 12
$ws.Cells.Item(12, 3).Value = 15
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(14, 3).Value = 5
$ws.Cells.Item(16, 3).Value = 11
$ws.Cells.Item(19, 3).Value = 13
$ws.Cells.Item(20, 3).Value = 14
$ws.Cells.Item(21, 3).Value = 27
$ws.Cells.Item(22, 3).Value = 13
$ws.Cells.Item(24, 3).Value = 22
$ws.Cells.Item(27, 3).Value = 19
$ws.Cells.Item(28, 3).Value = 21
$ws.Cells.Item(34, 3).Value = 6
$ws.Cells.Item(36, 3).Value = 17
$ws.Cells.Item(37, 3).Value = 19
$ws.Cells.Item(40, 3).Value = 11
$ws.Cells.Item(42, 3).Value = 5
$ws.Cells.Item(43, 3).Value = 13
$ws.Cells.Item(45, 3).Value = 25
$ws.Cells.Item(46, 3).Value = 8
$ws.Cells.Item(50, 3).Value = 5
$ws.Cells.Item(51, 3).Value = 14
$ws.Cells.Item(53, 3).Value = 24
$ws.Cells.Item(55, 3).Value = 15
$ws.Cells.Item(59, 3).Value = 19
$ws.Cells.Item(61, 3).Value = 28
$ws.Cells.Item(63, 3).Value = 10
$ws.Cells.Item(65, 3).Value = 16
$ws.Cells.Item(68, 3).Value = 17
$ws.Cells.Item(71, 3).Value = 9
$ws.Cells.Item(72, 3).Value = 11
$ws.Cells.Item(73, 3).Value = 19
$ws.Cells.Item(76, 3).Value = 19
$ws.Cells.Item(77, 3).Value = 22
$ws.Cells.Item(78, 3).Value = 15
$ws.Cells.Item(79, 3).Value = 18
$ws.Cells.Item(80, 3).Value = 21
$ws.Cells.Item(81, 3).Value = 24
$ws.Cells.Item(82, 3).Value = 4
$ws.Cells.Item(88, 3).Value = 15
$ws.Cells.Item(95, 3).Value = 7
$ws.Cells.Item(97, 3).Value = 15
$ws.Cells.Item(104, 3).Value = 19
$ws.Cells.Item(105, 3).Value = 25
$ws.Cells.Item(107, 3).Value = 11
$ws.Cells.Item(108, 3).Value = 17
$ws.Cells.Item(109, 3).Value = 21
$ws.Cells.Item(111, 3).Value = 7
$ws.Cells.Item(114, 3).Value = 5
$ws.Cells.Item(115, 3).Value = 13
$ws.Cells.Item(116, 3).Value = 14
$ws.Cells.Item(117, 3).Value = 19
$ws.Cells.Item(120, 3).Value = 19
$ws.Cells.Item(121, 3).Value = 26
$ws.Cells.Item(122, 3).Value = 8
$ws.Cells.Item(127, 3).Value = 13
$ws.Cells.Item(130, 3).Value = 12
$ws.Cells.Item(132, 3).Value = 23
$ws.Cells.Item(134, 3).Value = 9
$ws.Cells.Item(135, 3).Value = 15
$ws.Cells.Item(136, 3).Value = 18
$ws.Cells.Item(137, 3).Value = 23
$ws.Cells.Item(139, 3).Value = 13
$ws.Cells.Item(140, 3).Value = 19
$ws.Cells.Item(143, 3).Value = 15
$ws.Cells.Item(145, 3).Value = 19
$ws.Cells.Item(148, 3).Value = 22
$ws.Cells.Item(152, 3).Value = 23
$ws.Cells.Item(154, 3).Value = 4
$ws.Cells.Item(155, 3).Value = 11
$ws.Cells.Item(156, 3).Value = 14
$ws.Cells.Item(157, 3).Value = 19
$ws.Cells.Item(159, 3).Value = 8
$ws.Cells.Item(161, 3).Value = 14
$ws.Cells.Item(162, 3).Value = 8
$ws.Cells.Item(163, 3).Value = 11
$ws.Cells.Item(164, 3).Value = 14
$ws.Cells.Item(165, 3).Value = 25
$ws.Cells.Item(173, 3).Value = 18
$ws.Cells.Item(175, 3).Value = 8
$ws.Cells.Item(176, 3).Value = 9
$ws.Cells.Item(177, 3).Value = 13
$ws.Cells.Item(180, 3).Value = 20
$ws.Cells.Item(182, 3).Value = 10
$ws.Cells.Item(183, 3).Value = 16
$ws.Cells.Item(184, 3).Value = 19
$ws.Cells.Item(185, 3).Value = 22
$ws.Cells.Item(186, 3).Value = 13
$ws.Cells.Item(188, 3).Value = 21
$ws.Cells.Item(191, 3).Value = 19
$ws.Cells.Item(199, 3).Value = 17
$ws.Cells.Item(203, 3).Value = 15
$ws.Cells.Item(208, 3).Value = 17
$ws.Cells.Item(209, 3).Value = 20
$ws.Cells.Item(212, 3).Value = 13
$ws.Cells.Item(213, 3).Value = 15
$ws.Cells.Item(216, 3).Value = 13
$ws.Cells.Item(217, 3).Value = 18
$ws.Cells.Item(218, 3).Value = 6
$ws.Cells.Item(221, 3).Value = 14
$ws.Cells.Item(225, 3).Value = 22
$ws.Cells.Item(226, 3).Value = 9
$ws.Cells.Item(227, 3).Value = 15
$ws.Cells.Item(229, 3).Value = 20
